$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.213.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.03%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.855.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.36%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6982'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.79%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07760'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.45%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3067'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.10%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07809'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.40%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.861.48'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.24%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.097'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.35%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.51%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6858'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.39%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.531'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.64%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008432'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.54%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.214.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.05%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.29%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.111.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.92%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.08%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.12%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.502'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.18%  '

# Row 24
$ws.Range("E24").Value = '  -0.05%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1495'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.79%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.96%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.860'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.37%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.63%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.561'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.18%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.238'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.92%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.188'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.80%  '

# Row 32
$ws.Range("E32").Value = '  -1.08%  '

# Row 33
$ws.Range("E33").Value = '  -0.94%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7597'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.13%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.844'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.21%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.166'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.95%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.709'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.11%  '

# Row 38
$ws.Range("E38").Value = '  -0.17%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.211.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.42%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.722'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.66%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8987'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.05%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.64'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.08%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9992'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.12%  '

# Row 44
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.011.66'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.11%  '

# Row 45
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.516'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -12.10%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000125'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.71%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.56%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5177'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.44%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.534'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.39%  '

# Row 50
$ws.Range("E50").Value = '  -1.64%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.011'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.36%  '
